$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF") using same formatting as other headers (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-18
$data = @(
    @(9, 9),
    @(4, 5),
    @(8, 8),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(6, 6),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
